# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the last (58fa31bb...) row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-27 08:06:32"
$wsZhCn.Range("G4").Value = "2016-01-27 08:07:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-27 08:06:42"
$wsDeDe.Range("G4").Value = "2016-01-27 08:07:39"
